$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ C = 64.37254901960785; D = 66.44061302681992; E = 65.7235494880546;  F = 64.76140350877193; G = 65.1549815498155 }
    3 = @{ C = 65.68774703557312; D = 63.18992248062015;  E = 64.66283524904215; F = 64.87867647058823; G = 63.13454545454545 }
    4 = @{ C = 64.06209150326798; D = 65.18493150684931;  E = 65.70819672131148; F = 63.46229508196721; G = 66.77397260273973 }
    5 = @{ C = 66.45270270270271; D = 64.31660231660231;  E = 64.79180887372014; F = 64.77737226277372; G = 64.30584192439862 }
    6 = @{ C = 63.55063291139241; D = 65.17081850533808;  E = 64.90977443609023; F = 66.07986111111111; G = 65.06293706293707 }
    7 = @{ C = 63.8841059602649;  D = 65.67883211678833;  E = 66.32432432432432; F = 66.67870036101083; G = 64.55892255892256 }
    8 = @{ C = 62.31428571428572; D = 63.87713310580205;  E = 64.94809688581314; F = 67.38698630136986; G = 65.98961937716263 }
    9 = @{ C = 65.01730103806229; D = 63.24414715719063;  E = 64.29818181818182; F = 66.8301282051282;  G = 64.20265780730897 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
